$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "platforms" column (O) for several products: PC -> Steam / Software
$ws.Range("O4").Value = "Steam,PlayStation 5,Xbox Series X"   # Helldrivers 2
$ws.Range("O5").Value = "Steam,PlayStation 5"   # Nioh 3
$ws.Range("O6").Value = "Software"               # Office 365
$ws.Range("O7").Value = "Steam"                  # Terraria
$ws.Range("O8").Value = "Steam,PlayStation 5"   # The Witcher 3: Wild Hunt

# Move the viewport / active selection as recorded by the editing session
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("J17").Select()
